$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = 44441
$ws.Range("D3").Value = 44441
$ws.Range("D2:D3").NumberFormat = "mm-dd-yy"
